# 02/09 avance prog concilia auto 2h
#
# Insert a new row above row 5 (pushes the existing "Comercio" mini-table and
# everything below it down by one row), then fill the newly freed row 6 with
# a 0-10 numbering header, and finally leave the selection on E23 to match
# where work left off.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert one entire row at row 5 - shifts rows 5..35 down to 6..36.
$ws.Rows("5:5").Insert()

# Populate the newly inserted (and now blank) row 6 with a 0-10 sequence.
$ws.Range("C6").Value = 0
$ws.Range("D6").Value = 1
$ws.Range("E6").Value = 2
$ws.Range("F6").Value = 3
$ws.Range("G6").Value = 4
$ws.Range("H6").Value = 5
$ws.Range("I6").Value = 6
$ws.Range("J6").Value = 7
$ws.Range("K6").Value = 8
$ws.Range("L6").Value = 9
$ws.Range("M6").Value = 10

# Leave the selection where the user last was when saving.
$ws.Range("E23").Select() | Out-Null
